{"js": "const replacements = [\n  [\"2023-02-18 Saturday\", \"2023-02-19 Sunday\"],\n  [\"75-35=\", \"35-31=\"],\n  [\"27+70=\", \"36-12=\"],\n  [\"48-10=\", \"5+67=\"],\n  [\"37+14=\", \"77+2=\"],\n  [\"95-55=\", \"46-30=\"],\n  [\"56+16=\", \"41+42=\"],\n  [\"9+42=\", \"73-64=\"],\n  [\"75-47=\", \"6+83=\"],\n  [\"60+8=\", \"5+56=\"],\n  [\"67+6=\", \"31+65=\"],\n  [\"65+7=\", \"60+4=\"],\n  [\"20+23=\", \"5+82=\"],\n  [\"70-47=\", \"17+31=\"],\n  [\"73-43=\", \"15+34=\"],\n  [\"46+5=\", \"14+69=\"],\n  [\"91-84=\", \"18+55=\"],\n  [\"46+38=\", \"46+29=\"],\n  [\"41-28=\", \"5+25=\"],\n  [\"31+48=\", \"13+40=\"],\n  [\"85-48=\", \"23+37=\"],\n  [\"85-26=\", \"0+1=\"],\n  [\"74-17=\", \"64-29=\"],\n  [\"14-9=\", \"81-74=\"],\n  [\"12+67=\", \"18-0=\"],\n  [\"10+1=\", \"82-40=\"],\n  [\"75-15=\", \"33+44=\"],\n  [\"60+22=\", \"2+31=\"],\n  [\"3+37=\", \"63-52=\"],\n  [\"91-56=\", \"6+81=\"],\n  [\"30-14=\", \"93+4=\"],\n  [\"21+1=\", \"43-6=\"],\n  [\"29-5=\", \"84+7=\"],\n  [\"98-24=\", \"47-9=\"],\n  [\"69-40=\", \"89-41=\"],\n  [\"35-4=\", \"68+26=\"],\n  [\"20+46=\", \"4+50=\"],\n  [\"57-48=\", \"96-61=\"],\n  [\"70-57=\", \"28-24=\"],\n  [\"92-11=\", \"90-70=\"],\n  [\"99-95=\", \"96-61=\"],\n  [\"69-11=\", \"93-17=\"],\n  [\"78-73=\", \"14+73=\"],\n  [\"75-53=\", \"62+7=\"],\n  [\"14+13=\", \"4+14=\"],\n  [\"64-1=\", \"98-74=\"],\n  [\"99-74=\", \"61+35=\"],\n  [\"39+28=\", \"11+6=\"],\n  [\"77-6=\", \"20+12=\"],\n  [\"93-6=\", \"84-31=\"],\n  [\"89-14=\", \"10+87=\"],\n  [\"87-26=\", \"65-4=\"],\n  [\"65+24=\", \"73-4=\"],\n  [\"43+34=\", \"67-23=\"],\n  [\"58+4=\", \"46-20=\"],\n  [\"48+51=\", \"55+20=\"],\n  [\"35-25=\", \"47+10=\"],\n  [\"66-6=\", \"46-25=\"],\n  [\"75-31=\", \"65+9=\"],\n  [\"32+62=\", \"67-26=\"],\n  [\"69-28=\", \"11+63=\"],\n  [\"67+27=\", \"52+39=\"],\n  [\"11+66=\", \"37+34=\"],\n  [\"73+2=\", \"72-19=\"],\n  [\"75-74=\", \"30-15=\"],\n  [\"30-4=\", \"26+43=\"],\n  [\"36+44=\", \"49-15=\"],\n  [\"26-5=\", \"15+46=\"],\n  [\"89-30=\", \"69-24=\"],\n  [\"47-29=\", \"22+70=\"],\n  [\"69-39=\", \"20+29=\"],\n  [\"75+10=\", \"21+48=\"],\n  [\"5+70=\", \"65+4=\"],\n  [\"95+0=\", \"36+12=\"],\n  [\"51+34=\", \"82-27=\"],\n  [\"40+15=\", \"96-35=\"],\n  [\"86-43=\", \"98-34=\"],\n  [\"53-24=\", \"0+56=\"],\n  [\"6+39=\", \"91+8=\"],\n  [\"38+18=\", \"25+22=\"],\n  [\"71-23=\", \"4+95=\"],\n  [\"41-6=\", \"43+32=\"],\n  [\"88-80=\", \"49-32=\"],\n  [\"59-26=\", \"87-2=\"],\n  [\"75-21=\", \"43+3=\"],\n  [\"74+19=\", \"45-32=\"],\n  [\"44-26=\", \"68-55=\"],\n  [\"71-50=\", \"84+10=\"],\n  [\"95-69=\", \"88-31=\"],\n  [\"1+46=\", \"27+15=\"],\n  [\"43+52=\", \"55+5=\"],\n  [\"71+14=\", \"16-3=\"],\n  [\"46+24=\", \"53-50=\"],\n  [\"55+3=\", \"33+21=\"],\n  [\"36+11=\", \"80-3=\"],\n  [\"71+15=\", \"41+34=\"],\n  [\"1+67=\", \"86-29=\"],\n  [\"96-68=\", \"3+83=\"],\n  [\"20+78=\", \"26+59=\"],\n  [\"70+14=\", \"98-41=\"],\n  [\"99-10=\", \"98-28=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"2023-02-18 Saturday\"; After = \"2023-02-19 Sunday\" },\n    @{ Before = \"75-35=\"; After = \"35-31=\" },\n    @{ Before = \"27+70=\"; After = \"36-12=\" },\n    @{ Before = \"48-10=\"; After = \"5+67=\" },\n    @{ Before = \"37+14=\"; After = \"77+2=\" },\n    @{ Before = \"95-55=\"; After = \"46-30=\" },\n    @{ Before = \"56+16=\"; After = \"41+42=\" },\n    @{ Before = \"9+42=\"; After = \"73-64=\" },\n    @{ Before = \"75-47=\"; After = \"6+83=\" },\n    @{ Before = \"60+8=\"; After = \"5+56=\" },\n    @{ Before = \"67+6=\"; After = \"31+65=\" },\n    @{ Before = \"65+7=\"; After = \"60+4=\" },\n    @{ Before = \"20+23=\"; After = \"5+82=\" },\n    @{ Before = \"70-47=\"; After = \"17+31=\" },\n    @{ Before = \"73-43=\"; After = \"15+34=\" },\n    @{ Before = \"46+5=\"; After = \"14+69=\" },\n    @{ Before = \"91-84=\"; After = \"18+55=\" },\n    @{ Before = \"46+38=\"; After = \"46+29=\" },\n    @{ Before = \"41-28=\"; After = \"5+25=\" },\n    @{ Before = \"31+48=\"; After = \"13+40=\" },\n    @{ Before = \"85-48=\"; After = \"23+37=\" },\n    @{ Before = \"85-26=\"; After = \"0+1=\" },\n    @{ Before = \"74-17=\"; After = \"64-29=\" },\n    @{ Before = \"14-9=\"; After = \"81-74=\" },\n    @{ Before = \"12+67=\"; After = \"18-0=\" },\n    @{ Before = \"10+1=\"; After = \"82-40=\" },\n    @{ Before = \"75-15=\"; After = \"33+44=\" },\n    @{ Before = \"60+22=\"; After = \"2+31=\" },\n    @{ Before = \"3+37=\"; After = \"63-52=\" },\n    @{ Before = \"91-56=\"; After = \"6+81=\" },\n    @{ Before = \"30-14=\"; After = \"93+4=\" },\n    @{ Before = \"21+1=\"; After = \"43-6=\" },\n    @{ Before = \"29-5=\"; After = \"84+7=\" },\n    @{ Before = \"98-24=\"; After = \"47-9=\" },\n    @{ Before = \"69-40=\"; After = \"89-41=\" },\n    @{ Before = \"35-4=\"; After = \"68+26=\" },\n    @{ Before = \"20+46=\"; After = \"4+50=\" },\n    @{ Before = \"57-48=\"; After = \"96-61=\" },\n    @{ Before = \"70-57=\"; After = \"28-24=\" },\n    @{ Before = \"92-11=\"; After = \"90-70=\" },\n    @{ Before = \"99-95=\"; After = \"96-61=\" },\n    @{ Before = \"69-11=\"; After = \"93-17=\" },\n    @{ Before = \"78-73=\"; After = \"14+73=\" },\n    @{ Before = \"75-53=\"; After = \"62+7=\" },\n    @{ Before = \"14+13=\"; After = \"4+14=\" },\n    @{ Before = \"64-1=\"; After = \"98-74=\" },\n    @{ Before = \"99-74=\"; After = \"61+35=\" },\n    @{ Before = \"39+28=\"; After = \"11+6=\" },\n    @{ Before = \"77-6=\"; After = \"20+12=\" },\n    @{ Before = \"93-6=\"; After = \"84-31=\" },\n    @{ Before = \"89-14=\"; After = \"10+87=\" },\n    @{ Before = \"87-26=\"; After = \"65-4=\" },\n    @{ Before = \"65+24=\"; After = \"73-4=\" },\n    @{ Before = \"43+34=\"; After = \"67-23=\" },\n    @{ Before = \"58+4=\"; After = \"46-20=\" },\n    @{ Before = \"48+51=\"; After = \"55+20=\" },\n    @{ Before = \"35-25=\"; After = \"47+10=\" },\n    @{ Before = \"66-6=\"; After = \"46-25=\" },\n    @{ Before = \"75-31=\"; After = \"65+9=\" },\n    @{ Before = \"32+62=\"; After = \"67-26=\" },\n    @{ Before = \"69-28=\"; After = \"11+63=\" },\n    @{ Before = \"67+27=\"; After = \"52+39=\" },\n    @{ Before = \"11+66=\"; After = \"37+34=\" },\n    @{ Before = \"73+2=\"; After = \"72-19=\" },\n    @{ Before = \"75-74=\"; After = \"30-15=\" },\n    @{ Before = \"30-4=\"; After = \"26+43=\" },\n    @{ Before = \"36+44=\"; After = \"49-15=\" },\n    @{ Before = \"26-5=\"; After = \"15+46=\" },\n    @{ Before = \"89-30=\"; After = \"69-24=\" },\n    @{ Before = \"47-29=\"; After = \"22+70=\" },\n    @{ Before = \"69-39=\"; After = \"20+29=\" },\n    @{ Before = \"75+10=\"; After = \"21+48=\" },\n    @{ Before = \"5+70=\"; After = \"65+4=\" },\n    @{ Before = \"95+0=\"; After = \"36+12=\" },\n    @{ Before = \"51+34=\"; After = \"82-27=\" },\n    @{ Before = \"40+15=\"; After = \"96-35=\" },\n    @{ Before = \"86-43=\"; After = \"98-34=\" },\n    @{ Before = \"53-24=\"; After = \"0+56=\" },\n    @{ Before = \"6+39=\"; After = \"91+8=\" },\n    @{ Before = \"38+18=\"; After = \"25+22=\" },\n    @{ Before = \"71-23=\"; After = \"4+95=\" },\n    @{ Before = \"41-6=\"; After = \"43+32=\" },\n    @{ Before = \"88-80=\"; After = \"49-32=\" },\n    @{ Before = \"59-26=\"; After = \"87-2=\" },\n    @{ Before = \"75-21=\"; After = \"43+3=\" },\n    @{ Before = \"74+19=\"; After = \"45-32=\" },\n    @{ Before = \"44-26=\"; After = \"68-55=\" },\n    @{ Before = \"71-50=\"; After = \"84+10=\" },\n    @{ Before = \"95-69=\"; After = \"88-31=\" },\n    @{ Before = \"1+46=\"; After = \"27+15=\" },\n    @{ Before = \"43+52=\"; After = \"55+5=\" },\n    @{ Before = \"71+14=\"; After = \"16-3=\" },\n    @{ Before = \"46+24=\"; After = \"53-50=\" },\n    @{ Before = \"55+3=\"; After = \"33+21=\" },\n    @{ Before = \"36+11=\"; After = \"80-3=\" },\n    @{ Before = \"71+15=\"; After = \"41+34=\" },\n    @{ Before = \"1+67=\"; After = \"86-29=\" },\n    @{ Before = \"96-68=\"; After = \"3+83=\" },\n    @{ Before = \"20+78=\"; After = \"26+59=\" },\n    @{ Before = \"70+14=\"; After = \"98-41=\" },\n    @{ Before = \"99-10=\"; After = \"98-28=\" },\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $pair.Before,  # FindText\n        $true,         # MatchCase\n        $true,         # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $pair.After,   # ReplaceWith\n        2              # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
